# edit.ps1 - applies the "_模板.xlsx" template-table revision:
#  - renames the example sheets to the new alias/subtable/split-table naming scheme
#  - updates the explanatory text cells on "title前六行" and "解释" to describe
#    the new separators (| for list, & for map, + for subtable, = for alias)
#  - adds a new "*****" marker cell and explanatory "eg:" example cells

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheets (tab order / sheetId / rId are untouched - only the name).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(3).Name = '=rare#别名例子'
$wb.Worksheets.Item(4).Name = '+subList#list子表'
$wb.Worksheets.Item(5).Name = '+subMap#map子表'
$wb.Worksheets.Item(6).Name = '总表A'
$wb.Worksheets.Item(7).Name = '总表A-分表1#分表例子'
$wb.Worksheets.Item(8).Name = '总表A-分表2#分表例子'

# ---------------------------------------------------------------------------
# 2) "title前六行" sheet: update the relation/alias explanation cell.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item('title前六行')
$ws1.Range('B8').Value = '关联(别名/子表/分表:参考页签:解释)'

# ---------------------------------------------------------------------------
# 3) "解释" sheet: update existing explanation text + add new "eg:" example
#    cells in column E (styled like the existing "差" / Bad cells).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item('解释')

# row 10 - split-table explanation
$ws2.Range('D10').Value = '同一个excel下的sheet, sheet名为aaa-bbb的表会被归为aaa的分表'
$ws2.Range('E10').Value = 'eg:同一个excel下的sheet, "道具-1","道具-2" 都合并到"道具"'
$ws2.Range('E10').Style = '差'

# row 12 - list separator explanation (";" -> "|")
$ws2.Range('D12').Value = '分割符号为英文下的|'
$ws2.Range('E12').Value = 'eg: 1|2|3|4'
$ws2.Range('E12').Style = '差'

# row 13 - map separator explanation (";" -> "&")
$ws2.Range('D13').Value = 'k和v的分割符号为|,多个kv之间的分割为&'
$ws2.Range('E13').Value = 'eg: 1|a&1|b&2|c'
$ws2.Range('E13').Style = '差'

# row 14 - value unchanged, but E14 now carries the "差" style (was default).
$ws2.Range('E14').Style = '差'

# row 16 - list subtable separator explanation ("|" -> "+")
$ws2.Range('D16').Value = 'title的第4行,格式为:表名+Sheet名字'
$ws2.Range('E16').Value = '注:表名可省略,表示本表,但"+"不可省略'

# row 17 - map subtable separator explanation ("|" -> "+")
$ws2.Range('D17').Value = 'title的第4行,格式为:表名+Sheet名字'
$ws2.Range('E17').Value = '注:表名可省略,表示本表,但"+"不可省略'

# row 18 - alias separator explanation ("." -> "=")
$ws2.Range('D18').Value = 'title的第4行,格式为:表名=Sheet名字'
$ws2.Range('E18').Value = '注:表名可省略,表示本表,但"="不可省略'

# page setup: print as portrait on A4-ish paper (paperSize 9)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4) "=rare#别名例子" sheet (was "别名例子"): add the new "*****" marker cell.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item('=rare#别名例子')
$ws3.Range('K30').Value = '*****'
